# Tutorial 6 - updated solution
# Updates the Date column (A) to use dashes instead of slashes, and
# refreshes the Real/Duplicate/Invalid/Absent derived counts in columns
# D-H for the rows whose attendance status changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write the date text into column A without letting Excel's
# autodetection reinterpret ambiguous "dd-mm-yyyy" strings (where
# dd <= 12) as a date serial number - it must stay plain text, exactly
# like the other (unambiguous) rows already do.
function Set-DateText($row, $text) {
    $cell = $ws.Range("A$row")
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# r  Date           D(Total) E(Real) G(Invalid) H(Absent)
Set-DateText 3 "28-07-2022"
$ws.Range("D3").Value  = 1
$ws.Range("G3").Value  = 1

Set-DateText 4 "01-08-2022"
$ws.Range("D4").Value  = 1
$ws.Range("E4").Value  = 1
$ws.Range("H4").Value  = 0

Set-DateText 5 "04-08-2022"
$ws.Range("D5").Value  = 1
$ws.Range("E5").Value  = 1
$ws.Range("H5").Value  = 0

Set-DateText 6 "08-08-2022"

Set-DateText 7 "11-08-2022"

Set-DateText 8 "15-08-2022"

Set-DateText 9 "18-08-2022"

Set-DateText 10 "22-08-2022"
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1
$ws.Range("H10").Value = 0

Set-DateText 11 "25-08-2022"
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 1
$ws.Range("H11").Value = 0

Set-DateText 12 "29-08-2022"
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 1
$ws.Range("H12").Value = 0

Set-DateText 13 "01-09-2022"

Set-DateText 14 "05-09-2022"
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 1
$ws.Range("H14").Value = 0

Set-DateText 15 "08-09-2022"

Set-DateText 16 "12-09-2022"

Set-DateText 17 "15-09-2022"

Set-DateText 18 "19-09-2022"

Set-DateText 19 "22-09-2022"

Set-DateText 20 "26-09-2022"
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 1
$ws.Range("H20").Value = 0

Set-DateText 21 "29-09-2022"
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 1
$ws.Range("H21").Value = 0
